# POC of inventory search - add a fourth "Implemented" column and three
# new questions to the Questions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D values (boolean) for existing rows 2-11 ---
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $true
$ws.Range("D4").Value = $true
$ws.Range("D5").Value = $false
$ws.Range("D6").Value = $false
$ws.Range("D7").Value = $true
$ws.Range("D8").Value = $true
$ws.Range("D9").Value = $false
$ws.Range("D10").Value = $true
$ws.Range("D11").Value = $true

# --- Three new questions appended as rows 12-14 ---
$ws.Range("A12").Value = "What's the total value of my stock?"
$ws.Range("D12").Value = $false

$ws.Range("A13").Value = "Which stone in my inventory has least price?"
$ws.Range("D13").Value = $false

$ws.Range("A14").Value = "Which category has most number of stones?"
$ws.Range("D14").Value = $false

# Rows 12 and 13 in column D pick up the "Bad" highlight style used
# elsewhere in the sheet (same as B7) to flag them as not yet answered.
$ws.Range("D12").Style = "Bad"
$ws.Range("D13").Style = "Bad"

# Move the active selection to A4, matching the author's saved cursor
# position when they committed this change.
$ws.Range("A4").Select()
